$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "26.013.83"
$ws.Range("E2").Value = "  +0.61%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.762.04"
$ws.Range("E3").Value = "  +0.54%  "
$ws.Range("E4").Value = "  +0.38%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "237.84"
$ws.Range("E5").Value = "  -0.78%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.002"
$ws.Range("E6").Value = "  +0.28%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5227"
$ws.Range("E7").Value = "  +2.30%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2755"
$ws.Range("E8").Value = "  -1.21%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "40.56"
$ws.Range("E9").Value = "  -3.93%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.06199"
$ws.Range("E10").Value = "  -0.16%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.775.07"
$ws.Range("E11").Value = "  +1.26%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.07021"
$ws.Range("E12").Value = "  +0.74%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "15.67"
$ws.Range("E13").Value = "  -0.86%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.6388"
$ws.Range("E14").Value = "  +4.72%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "4.547"
$ws.Range("E15").Value = "  +0.53%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "78.16"
$ws.Range("E16").Value = "  +0.58%  "
$ws.Range("E17").Value = "  +0.44%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "1.001"
$ws.Range("E18").Value = "  +0.26%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "26.048.41"
$ws.Range("E19").Value = "  +0.78%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "11.64"
$ws.Range("E20").Value = "  -0.60%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.000006731"
$ws.Range("E21").Value = "  -3.75%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "2.001.79"
$ws.Range("E22").Value = "  +1.60%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.075"
$ws.Range("E23").Value = "  -0.33%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "8.455"
$ws.Range("E24").Value = "  +2.77%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "5.197"
$ws.Range("E25").Value = "  -1.27%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "138.85"
$ws.Range("E26").Value = "  +0.76%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.512"
$ws.Range("E27").Value = "  +2.71%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.842"
$ws.Range("E28").Value = "  +0.69%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "15.17"
$ws.Range("E29").Value = "  +0.67%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "103.00"
$ws.Range("E30").Value = "  -0.83%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.08385"
$ws.Range("E31").Value = "  +2.36%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.694"
$ws.Range("E32").Value = "  -0.31%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.444"
$ws.Range("E33").Value = "  -2.54%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.04464"
$ws.Range("E34").Value = "  -1.60%  "
$ws.Range("E35").Value = "  -0.10%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.001"
$ws.Range("E36").Value = "  +0.65%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.6058"
$ws.Range("E37").Value = "  -1.23%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.741"
$ws.Range("E38").Value = "  +1.66%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01588"
$ws.Range("E39").Value = "  +1.91%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.975"
$ws.Range("E40").Value = "  +3.16%  "
$ws.Range("E41").Value = "  +0.27%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "102.72"
$ws.Range("E42").Value = "  -0.70%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.3874"
$ws.Range("E43").Value = "  -0.36%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.7445"
$ws.Range("E44").Value = "  -0.59%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "4.931"
$ws.Range("E45").Value = "  -0.29%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.05516"
$ws.Range("E46").Value = "  +2.12%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "6.338"
$ws.Range("E47").Value = "  +5.01%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.1118"
$ws.Range("E48").Value = "  +0.37%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "30.17"
$ws.Range("E49").Value = "  -0.25%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "52.58"
$ws.Range("E50").Value = "  -0.65%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.004"
$ws.Range("E51").Value = "  +0.77%  "
